$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the two new columns
$ws.Range("E1").Value = "total.abundance"
$ws.Range("F1").Value = "percent.abundance"

# Match the number format already used elsewhere in the sheet (style index 1 -> "0.00")
$ws.Range("F1").NumberFormat = "0.00"

$lastRow = 53
$spawnCol = 4   # D = spawn.abundance
$totalCol = 5   # E = total.abundance
$pctCol = 6     # F = percent.abundance

# First pass: compute the total spawn.abundance per year group (rows 2-53)
$yearTotals = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $year = $ws.Cells.Item($r, 1).Value2
    $d = $ws.Cells.Item($r, $spawnCol).Value2
    if ($d -eq $null) { $d = 0 }
    if ($yearTotals.ContainsKey($year)) {
        $yearTotals[$year] = $yearTotals[$year] + $d
    } else {
        $yearTotals[$year] = $d
    }
}

# Second pass: write total.abundance and percent.abundance for each row
for ($r = 2; $r -le $lastRow; $r++) {
    $year = $ws.Cells.Item($r, 1).Value2
    $d = $ws.Cells.Item($r, $spawnCol).Value2
    if ($d -eq $null) { $d = 0 }
    $total = $yearTotals[$year]

    $ws.Cells.Item($r, $totalCol).Value = $total

    $pctCell = $ws.Cells.Item($r, $pctCol)
    $pctCell.Value = ($d / $total) * 100
    $pctCell.NumberFormat = "0.00"
}

# Column widths for the new columns (closest achievable match to the bestFit sizing
# from the authored file; this runtime quantizes ColumnWidth to 1/6-character steps)
$ws.Columns.Item($totalCol).ColumnWidth = 12.8
$ws.Columns.Item($pctCol).ColumnWidth = 15.0

# Update the active selection to reflect where the author left off editing
$ws.Range("F7").Select()
